$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 158.1579
$ws.Range("I33").Value = 107.5
$ws.Range("J33").Value = 428.33334
$ws.Range("K33").Value = 107.5
$ws.Range("L33").Value = 428.33334
$ws.Range("M33").Value = 121.5
$ws.Range("N33").Value = -886.33334
# Row 133
$ws.Range("H133").Value = 12371.25
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 12371.25
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 12371.25
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -22491.25
# Row 137
$ws.Range("H137").Value = 17242340
$ws.Range("I137").Value = 22222842
$ws.Range("J137").Value = 2147.5386
$ws.Range("K137").Value = 66668526
$ws.Range("L137").Value = 6442.6158
$ws.Range("M137").Value = -66665976
$ws.Range("N137").Value = -11542.6158

$ws = $wb.Worksheets.Item("ARM")
# Row 30
$ws.Range("H30").Value = 2340000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 2340000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 2340000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -2340300
# Row 45
$ws.Range("H45").Value = 919.125
$ws.Range("I45").Value = 870.95
$ws.Range("J45").Value = 1160
$ws.Range("K45").Value = 870.95
$ws.Range("L45").Value = 1160
$ws.Range("M45").Value = -493.95
$ws.Range("N45").Value = -1914
# Row 74
$ws.Range("H74").Value = 4230.405
$ws.Range("I74").Value = 1124.1613
$ws.Range("J74").Value = 12984.363
$ws.Range("K74").Value = 1124.1613
$ws.Range("L74").Value = 12984.363
$ws.Range("M74").Value = -250.1613
$ws.Range("N74").Value = -14732.363
# Row 77
$ws.Range("H77").Value = 4230.405
$ws.Range("I77").Value = 1124.1613
$ws.Range("J77").Value = 12984.363
$ws.Range("K77").Value = 5620.8065
$ws.Range("L77").Value = 64921.815
$ws.Range("M77").Value = -1252.8065
$ws.Range("N77").Value = -73657.815
# Row 133
$ws.Range("H133").Value = 53399.832
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 53399.832
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 53399.832
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -58459.832
# Row 139
$ws.Range("H139").Value = 39517.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 39517.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 39517.75
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -49797.75

$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 49998.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 49998.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 49998.5
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -51692.5
# Row 105
$ws.Range("H105").Value = 297160.7
$ws.Range("I105").Value = 2823.0435
$ws.Range("J105").Value = 912593.9399999999
$ws.Range("K105").Value = 2823.0435
$ws.Range("L105").Value = 912593.9399999999
$ws.Range("M105").Value = -1076.0435
$ws.Range("N105").Value = -916087.9399999999

$ws = $wb.Worksheets.Item("CRP")
# Row 11
$ws.Range("H11").Value = 13780
$ws.Range("I11").Value = 27000
$ws.Range("J11").Value = 4966.6665
$ws.Range("K11").Value = 27000
$ws.Range("L11").Value = 4966.6665
$ws.Range("M11").Value = -26860
$ws.Range("N11").Value = -5246.6665
# Row 26
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("M26").Value = 10000
$ws.Range("N26").Value = -10574
# Row 31
$ws.Range("H31").Value = 1791.7709
$ws.Range("I31").Value = 1107.5526
$ws.Range("J31").Value = 4391.8
$ws.Range("K31").Value = 1107.5526
$ws.Range("L31").Value = 4391.8
$ws.Range("M31").Value = -812.5526
$ws.Range("N31").Value = -4981.8
# Row 34
$ws.Range("H34").Value = 1791.7709
$ws.Range("I34").Value = 1107.5526
$ws.Range("J34").Value = 4391.8
$ws.Range("K34").Value = 1107.5526
$ws.Range("L34").Value = 4391.8
$ws.Range("M34").Value = -905.5526
$ws.Range("N34").Value = -4795.8
# Row 58
$ws.Range("H58").Value = 2651.4285
$ws.Range("I58").Value = 1891.2
$ws.Range("J58").Value = 4552
$ws.Range("K58").Value = 1891.2
$ws.Range("L58").Value = 4552
$ws.Range("M58").Value = -1688.2
$ws.Range("N58").Value = -4958
# Row 105
$ws.Range("H105").Value = 875
$ws.Range("I105").Value = 875
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 875
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = 872
# Row 132
$ws.Range("H132").Value = 2785.4119
$ws.Range("I132").Value = 2409.1035
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 7227.310500000001
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -4697.310500000001
$ws.Range("N132").Value = -19964
# Row 136
$ws.Range("H136").Value = 2651.4285
$ws.Range("I136").Value = 1891.2
$ws.Range("J136").Value = 4552
$ws.Range("K136").Value = 5673.6
$ws.Range("L136").Value = 13656
$ws.Range("M136").Value = -3123.6
$ws.Range("N136").Value = -18756

$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 3167
$ws.Range("I70").Value = 2111.6667
$ws.Range("J70").Value = 4750
$ws.Range("K70").Value = 6335.000100000001
$ws.Range("L70").Value = 14250
$ws.Range("M70").Value = -6020.000100000001
$ws.Range("N70").Value = -14880
# Row 73
$ws.Range("H73").Value = 3167
$ws.Range("I73").Value = 2111.6667
$ws.Range("J73").Value = 4750
$ws.Range("K73").Value = 6335.000100000001
$ws.Range("L73").Value = 14250
$ws.Range("M73").Value = -5243.000100000001
$ws.Range("N73").Value = -16434
# Row 75
$ws.Range("H75").Value = 830.1
$ws.Range("I75").Value = 606.5
$ws.Range("J75").Value = 886
$ws.Range("K75").Value = 1819.5
$ws.Range("L75").Value = 2658
$ws.Range("M75").Value = -821.5
$ws.Range("N75").Value = -4654
# Row 78
$ws.Range("H78").Value = 830.1
$ws.Range("I78").Value = 606.5
$ws.Range("J78").Value = 886
$ws.Range("K78").Value = 5458.5
$ws.Range("L78").Value = 7974
$ws.Range("M78").Value = -466.5
$ws.Range("N78").Value = -17958
# Row 98
$ws.Range("H98").Value = 312.42856
$ws.Range("I98").Value = 199.6
$ws.Range("J98").Value = 594.5
$ws.Range("K98").Value = 598.8
$ws.Range("L98").Value = 1783.5
$ws.Range("M98").Value = 899.2
$ws.Range("N98").Value = -4779.5
# Row 105
$ws.Range("H105").Value = 6400
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 6400
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 19200
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -24442
# Row 124
$ws.Range("H124").Value = 1900
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 1900
$ws.Range("K124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("M124").Value = 5700
$ws.Range("N124").Value = -15520
# Row 125
$ws.Range("H125").Value = 2524.7083
$ws.Range("I125").Value = 1426.6666
$ws.Range("J125").Value = 2890.7222
$ws.Range("K125").Value = 4279.9998
$ws.Range("L125").Value = 8672.1666
$ws.Range("M125").Value = 640.0002000000004
$ws.Range("N125").Value = -18512.1666
# Row 131
$ws.Range("H131").Value = 7577602
$ws.Range("I131").Value = 578.625
$ws.Range("J131").Value = 9261385
$ws.Range("K131").Value = 1735.875
$ws.Range("L131").Value = 27784155
$ws.Range("M131").Value = 3304.125
$ws.Range("N131").Value = -27794235
# Row 140
$ws.Range("H140").Value = 6858.763
$ws.Range("I140").Value = 8490.111000000001
$ws.Range("J140").Value = 2854.5454
$ws.Range("K140").Value = 25470.333
$ws.Range("L140").Value = 8563.636200000001
$ws.Range("M140").Value = -20290.333
$ws.Range("N140").Value = -18923.6362

$ws = $wb.Worksheets.Item("GSM")
# Row 23
$ws.Range("H23").Value = 2937.5
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 6166.6665
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 6166.6665
$ws.Range("M23").Value = -777
$ws.Range("N23").Value = -6612.6665
# Row 138
$ws.Range("H138").Value = 73250
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 73250
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 73250
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -83530
# Row 139
$ws.Range("H139").Value = 59999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 59999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 59999
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -70279

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2811.2273
$ws.Range("I7").Value = 1609.4
$ws.Range("J7").Value = 3164.7058
$ws.Range("K7").Value = 1609.4
$ws.Range("L7").Value = 3164.7058
$ws.Range("M7").Value = -1497.4
$ws.Range("N7").Value = -3388.7058
# Row 126
$ws.Range("H126").Value = 2811.2273
$ws.Range("I126").Value = 1609.4
$ws.Range("J126").Value = 3164.7058
$ws.Range("K126").Value = 4828.200000000001
$ws.Range("L126").Value = 9494.117400000001
$ws.Range("M126").Value = -2358.200000000001
$ws.Range("N126").Value = -14434.1174
# Row 132
$ws.Range("H132").Value = 2603.6604
$ws.Range("I132").Value = 1738.6052
$ws.Range("J132").Value = 4795.1333
$ws.Range("K132").Value = 5215.8156
$ws.Range("L132").Value = 14385.3999
$ws.Range("M132").Value = -2685.8156
$ws.Range("N132").Value = -19445.3999

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2559.8193
$ws.Range("I132").Value = 2592.8772
$ws.Range("J132").Value = 2434.2
$ws.Range("K132").Value = 7778.6316
$ws.Range("L132").Value = 7302.599999999999
$ws.Range("M132").Value = -5248.6316
$ws.Range("N132").Value = -12362.6
# Row 136
$ws.Range("H136").Value = 34035.324
$ws.Range("I136").Value = 63600.625
$ws.Range("J136").Value = 2499
$ws.Range("K136").Value = 190801.875
$ws.Range("L136").Value = 7497
$ws.Range("M136").Value = -188251.875
$ws.Range("N136").Value = -12597
